$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the BI/PD/UM/MP symbol values between row 2 and row 3 (typo fix)
$ws.Range("AE2:AH2").Value = "I"
$ws.Range("AE3:AH3").Value = "A"

# Update the selected cell on the sheet view
$ws.Range("I13").Select()
